$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$replacements = @(
    @{Row=1;  Col=1; New="98×58="},
    @{Row=1;  Col=2; New="81×24="},
    @{Row=1;  Col=3; New="47×24="},
    @{Row=1;  Col=4; New="39×67="},
    @{Row=1;  Col=5; New="39×25="},

    @{Row=5;  Col=1; New="86×98="},
    @{Row=5;  Col=2; New="35×11="},
    @{Row=5;  Col=3; New="75×84="},
    @{Row=5;  Col=4; New="19×83="},
    @{Row=5;  Col=5; New="94×82="},

    @{Row=10; Col=1; New="71×40="},
    @{Row=10; Col=2; New="29×60="},
    @{Row=10; Col=3; New="87×26="},
    @{Row=10; Col=4; New="61×96="},
    @{Row=10; Col=5; New="86×53="},

    @{Row=15; Col=1; New="96×63="},
    @{Row=15; Col=2; New="18×39="},
    @{Row=15; Col=3; New="32×87="},
    @{Row=15; Col=4; New="90×88="},
    @{Row=15; Col=5; New="28×26="},

    @{Row=20; Col=1; New="52×57="},
    @{Row=20; Col=2; New="71×74="},
    @{Row=20; Col=3; New="20×68="},
    @{Row=20; Col=4; New="86×53="},
    @{Row=20; Col=5; New="70×95="}
)

foreach ($r in $replacements) {
    $cell = $t.Cell($r.Row, $r.Col)
    $cell.Range.Text = $r.New
}
